$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false | Out-Null

# ------------------------------------------------------------------
# 1. Delete the empty "Sheet1" worksheet (it was never populated and
#    is removed entirely from the workbook).
# ------------------------------------------------------------------
$wb.Worksheets.Item("Sheet1").Delete() | Out-Null

# ------------------------------------------------------------------
# 2. "Activity" sheet - add the two new Activity/Subject rows
#    (Meeting -> "Test Meeting", Other -> "Test Other").
# ------------------------------------------------------------------
$activityWs = $wb.Worksheets.Item("Activity")
$activityWs.Activate() | Out-Null
$activityWs.Range("B4").Value = "Test Meeting"
$activityWs.Range("B5").Value = "Test Other"
$activityWs.Range("B5").Select() | Out-Null

# ------------------------------------------------------------------
# 3. "Users" sheet - replace the sample user row with the new data.
# ------------------------------------------------------------------
$usersWs = $wb.Worksheets.Item("Users")
$usersWs.Activate() | Out-Null
$usersWs.Range("A1").Value = "CF Financial"
$usersWs.Range("A2").Value = "Thomas Bailey"
$usersWs.Range("B4").Select() | Out-Null

# ------------------------------------------------------------------
# 4. "ContactTypes" sheet - content is unchanged, just move the
#    cursor/selection as recorded in the saved file.
# ------------------------------------------------------------------
$contactTypesWs = $wb.Worksheets.Item("ContactTypes")
$contactTypesWs.Activate() | Out-Null
$contactTypesWs.Range("C21").Select() | Out-Null

# ------------------------------------------------------------------
# 5. "Contact" sheet - content is unchanged, becomes the active tab
#    with D20 selected.
# ------------------------------------------------------------------
$contactWs = $wb.Worksheets.Item("Contact")
$contactWs.Activate() | Out-Null
$contactWs.Range("D20").Select() | Out-Null
